$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RF)
$ws.Range("B3").Value = 0.261
$ws.Range("C3").Value = -0.066
$ws.Range("D3").Value = 0.469
$ws.Range("E3").Value = 0.6850000000000001
$ws.Range("F3").Value = 0.769
$ws.Range("G3").Value = 0.5570000000000001

# Row 4 (NN)
$ws.Range("B4").Value = 0.193
$ws.Range("C4").Value = -0.164
$ws.Range("D4").Value = 0.512
$ws.Range("E4").Value = 0.716
$ws.Range("F4").Value = 0.741
$ws.Range("G4").Value = 0.501

# Row 5 (RNN)
$ws.Range("B5").Value = 0.157
$ws.Range("C5").Value = -0.041
$ws.Range("D5").Value = 0.481
$ws.Range("E5").Value = 0.694
$ws.Range("F5").Value = 0.6870000000000001
$ws.Range("G5").Value = 0.555
